$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet view: selection moves from L9 to C7 (and the scrolled topLeftCell=A3 goes away) ---
[void]$ws.Range("C7").Select()

# --- A1: title - bump font size 15 -> 16, keep bold + Dark2 theme color, row height grows ---
$titleFont = $ws.Range("A1").Font
$titleFont.ThemeColor = 3   # msoThemeColorDark2 -> <color theme="3"/>
$titleFont.Size = 16
$ws.Rows(1).RowHeight = 21.75

# --- A2: replace the old "Project Introduction" rich text with the new "Business Scenario" text, now bold ---
$ws.Range("A2").Value = "Business Scenario:`nK & K JEANS is a retail clothing store specializing in branded jeans. The company places OEM orders and sells the products under its own brand name. Since orders must be placed before the start of each season, accurate demand forecasting and clear target setting are critical. Understanding the relationship between demand, order volume, and profitability will provide K & K JEANS with a competitive edge and help ensure sustainable growth."
$ws.Range("A2").Font.Bold = $true

# --- A4: replace the old "Objective" rich text with the new "Case Problem" text, now bold ---
$ws.Range("A4").Value = "Case Problem:  how to find out the break even point so not to over place order? "
$ws.Range("A4").Font.Bold = $true

# --- B11: Demand Q'ty changes from 4000 to 2475 (cascades through all dependent formulas) ---
$ws.Range("B11").Value = 2475

# --- C13: center the "FORMULAR" header label ---
$ws.Range("C13").HorizontalAlignment = -4108   # xlCenter

# --- E13: rename table header text ---
$ws.Range("E13").Value = "PROFIT TABLE BASED ON DEMAND AND ORDER QUANTITIES"

# --- New row 26: add the "Findings" paragraph under the profit table ---
$ws.Rows(26).RowHeight = 26.25
$ws.Range("A26").Value = "Findings:  the conditional formated table shows if the order quantity is 4500 pcs, the company will only start making profit when the demand quantity is around 2500 pcs, the break even point is  2475 pcs, at which the profit is zero."
$ws.Range("A26").Font.Bold = $true

# --- Force a pageSetup element to be written (portrait orientation) ---
$ws.PageSetup.Orientation = 1   # xlPortrait
